# ---------------------------------------------------------------------------
# valuations.xlsx -- "Fixed tests for fund_units"
#
# The sheet's sample/demo data (for a fund called "TSTF1") is replaced: new
# valuation dates, per-share values, portfolio-company names and instrument
# types. The old per-row hyperlinks (which pointed at a dev investment-
# instrument record) are removed, the used range grows out to column O (a
# side effect of pasting in a wider block of data from another workbook),
# and a stale AutoFilter database defined name is left behind referencing
# the new extent. Column D1 becomes the active/selected cell, and the sheet
# zoom is bumped to 113%.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Drop the six per-row hyperlinks that used to live in column E.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 2. Clear out the old data block (rows 2:7, cols A:E) so stray cells (e.g.
#    the always-empty column B) don't linger, then re-enter the new table.
# ---------------------------------------------------------------------------
$ws.Range("A2:E7").Clear()

# Row 2
$ws.Range("A2").Value = 45382
$ws.Range("C2").Value = 150
$ws.Range("D2").Value = "TSTF1 Port Co 1"
$ws.Range("E2").Value = "Equity"

# Row 3
$ws.Range("A3").Value = 45382
$ws.Range("C3").Value = 200
$ws.Range("D3").Value = "TSTF1 Port Co 1"
$ws.Range("E3").Value = "CCPS"

# Row 4
$ws.Range("A4").Value = 45382
$ws.Range("C4").Value = 200
$ws.Range("D4").Value = "TSTF1 Port Co 2"
$ws.Range("E4").Value = "CCPS"

# Row 5
$ws.Range("A5").Value = 45747
$ws.Range("C5").Value = 200
$ws.Range("D5").Value = "TSTF1 Port Co 1"
$ws.Range("E5").Value = "Equity"

# Row 6
$ws.Range("A6").Value = 45747
$ws.Range("C6").Value = 250
$ws.Range("D6").Value = "TSTF1 Port Co 1"
$ws.Range("E6").Value = "CCPS"

# Row 7
$ws.Range("A7").Value = 45747
$ws.Range("C7").Value = 250
$ws.Range("D7").Value = "TSTF1 Port Co 2"
$ws.Range("E7").Value = "CCPS"

# Date formatting for the Valuation Date column.
$ws.Range("A2:A7").NumberFormat = "m/d/yyyy"

# Restore the font/format used for the data rows (Arial 10).
$ws.Range("A2:E7").Font.Name = "Arial"
$ws.Range("A2:E7").Font.Size = 10

# ---------------------------------------------------------------------------
# 3. The pasted-in block extends the used range out to column O; columns
#    J:L on the lower rows carry a few blank formatted cells along for the
#    ride (left over from the source block's shape).
# ---------------------------------------------------------------------------
$ws.Range("J4:L7").Value = ""

# ---------------------------------------------------------------------------
# 4. A stale "_FilterDatabase" name (hidden, sheet-scoped) referencing the
#    new A1:O7 extent -- left behind from a filter applied on the source
#    data before it was pasted in.
# ---------------------------------------------------------------------------
$filterName = $ws.Names().Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$O`$7")
$filterName.Visible = $false

# ---------------------------------------------------------------------------
# 5. View state: active cell D1, zoom 113%.
# ---------------------------------------------------------------------------
$ws.Range("D1").Select()
$excel.ActiveWindow.Zoom = 113
